$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3: Test / Testowo / 03222222222 / test@XD.pl
$ws.Range("A3").Value = "Test"
$ws.Range("B3").Value = "Testowo"

# C3 looks like a pure number ("03222222222") - force text so the
# leading zero survives and it lands in sharedStrings like the rest.
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "03222222222"

$ws.Range("D3").Value = "test@XD.pl"

# Row 4: teststs322323 / 312123132 / 32132312321 / 321231312
$ws.Range("A4").Value = "teststs322323"

$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "312123132"

$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "32132312321"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "321231312"
